$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44425   # D2
$ws.Cells.Item(2, 10).Value = 400   # J2
$ws.Cells.Item(2, 11).Value = 11500   # K2
$ws.Cells.Item(2, 12).Value = 12000   # L2
$ws.Cells.Item(2, 13).Value = 11750   # M2
$ws.Cells.Item(2, 15).Value = "Provincia del Elquí"   # O2
$ws.Cells.Item(2, 16).Value = 470   # P2

$ws.Cells.Item(3, 4).Value = 44446   # D3
$ws.Cells.Item(3, 10).Value = 500   # J3
$ws.Cells.Item(3, 11).Value = 11000   # K3
$ws.Cells.Item(3, 12).Value = 12000   # L3
$ws.Cells.Item(3, 13).Value = 11500   # M3
$ws.Cells.Item(3, 15).Value = "Provincia del Elquí"   # O3
$ws.Cells.Item(3, 16).Value = 460   # P3

$ws.Cells.Item(4, 4).Value = 44694   # D4
$ws.Cells.Item(4, 10).Value = 480   # J4
$ws.Cells.Item(4, 11).Value = 17500   # K4
$ws.Cells.Item(4, 12).Value = 18000   # L4
$ws.Cells.Item(4, 13).Value = 17750   # M4
$ws.Cells.Item(4, 15).Value = "Provincia del Elquí"   # O4
$ws.Cells.Item(4, 16).Value = 710   # P4

$ws.Cells.Item(5, 4).Value = 44370   # D5
$ws.Cells.Item(5, 10).Value = 520   # J5
$ws.Cells.Item(5, 11).Value = 13000   # K5
$ws.Cells.Item(5, 12).Value = 14000   # L5
$ws.Cells.Item(5, 13).Value = 13500   # M5
$ws.Cells.Item(5, 15).Value = "Provincia del Elquí"   # O5
$ws.Cells.Item(5, 16).Value = 540   # P5

$ws.Cells.Item(6, 4).Value = 44721   # D6
$ws.Cells.Item(6, 10).Value = 500   # J6
$ws.Cells.Item(6, 11).Value = 14500   # K6
$ws.Cells.Item(6, 12).Value = 15000   # L6
$ws.Cells.Item(6, 13).Value = 14750   # M6
$ws.Cells.Item(6, 15).Value = "Provincia de Limarí"   # O6
$ws.Cells.Item(6, 16).Value = 590   # P6

$ws.Cells.Item(7, 4).Value = 44316   # D7
$ws.Cells.Item(7, 10).Value = 300   # J7
$ws.Cells.Item(7, 11).Value = 16000   # K7
$ws.Cells.Item(7, 12).Value = 17000   # L7
$ws.Cells.Item(7, 13).Value = 16500   # M7
$ws.Cells.Item(7, 15).Value = "Provincia del Elquí"   # O7
$ws.Cells.Item(7, 16).Value = 660   # P7

$ws.Cells.Item(8, 4).Value = 44376   # D8
$ws.Cells.Item(8, 10).Value = 400   # J8
$ws.Cells.Item(8, 11).Value = 12000   # K8
$ws.Cells.Item(8, 12).Value = 13000   # L8
$ws.Cells.Item(8, 13).Value = 12500   # M8
$ws.Cells.Item(8, 15).Value = "Provincia del Elquí"   # O8
$ws.Cells.Item(8, 16).Value = 500   # P8

$ws.Cells.Item(9, 4).Value = 44466   # D9
$ws.Cells.Item(9, 10).Value = 400   # J9
$ws.Cells.Item(9, 11).Value = 9500   # K9
$ws.Cells.Item(9, 12).Value = 10000   # L9
$ws.Cells.Item(9, 13).Value = 9750   # M9
$ws.Cells.Item(9, 15).Value = "Provincia del Elquí"   # O9
$ws.Cells.Item(9, 16).Value = 390   # P9

$ws.Cells.Item(10, 4).Value = 44377   # D10
$ws.Cells.Item(10, 10).Value = 520   # J10
$ws.Cells.Item(10, 11).Value = 12500   # K10
$ws.Cells.Item(10, 12).Value = 13000   # L10
$ws.Cells.Item(10, 13).Value = 12750   # M10
$ws.Cells.Item(10, 15).Value = "Provincia del Elquí"   # O10
$ws.Cells.Item(10, 16).Value = 510   # P10

$ws.Cells.Item(11, 4).Value = 44484   # D11
$ws.Cells.Item(11, 10).Value = 400   # J11
$ws.Cells.Item(11, 11).Value = 9000   # K11
$ws.Cells.Item(11, 12).Value = 10000   # L11
$ws.Cells.Item(11, 13).Value = 9500   # M11
$ws.Cells.Item(11, 15).Value = "Provincia del Elquí"   # O11
$ws.Cells.Item(11, 16).Value = 380   # P11

$ws.Cells.Item(12, 4).Value = 44714   # D12
$ws.Cells.Item(12, 10).Value = 400   # J12
$ws.Cells.Item(12, 11).Value = 14000   # K12
$ws.Cells.Item(12, 12).Value = 15000   # L12
$ws.Cells.Item(12, 13).Value = 14500   # M12
$ws.Cells.Item(12, 15).Value = "Provincia de Limarí"   # O12
$ws.Cells.Item(12, 16).Value = 580   # P12

$ws.Cells.Item(13, 4).Value = 44356   # D13
$ws.Cells.Item(13, 10).Value = 500   # J13
$ws.Cells.Item(13, 11).Value = 13000   # K13
$ws.Cells.Item(13, 12).Value = 14000   # L13
$ws.Cells.Item(13, 13).Value = 13500   # M13
$ws.Cells.Item(13, 15).Value = "Provincia de Limarí"   # O13
$ws.Cells.Item(13, 16).Value = 540   # P13

$ws.Cells.Item(14, 4).Value = 44372   # D14
$ws.Cells.Item(14, 10).Value = 500   # J14
$ws.Cells.Item(14, 11).Value = 13000   # K14
$ws.Cells.Item(14, 12).Value = 14000   # L14
$ws.Cells.Item(14, 13).Value = 13500   # M14
$ws.Cells.Item(14, 15).Value = "Provincia del Elquí"   # O14
$ws.Cells.Item(14, 16).Value = 540   # P14

$ws.Cells.Item(15, 4).Value = 44384   # D15
$ws.Cells.Item(15, 10).Value = 560   # J15
$ws.Cells.Item(15, 11).Value = 11500   # K15
$ws.Cells.Item(15, 12).Value = 12000   # L15
$ws.Cells.Item(15, 13).Value = 11750   # M15
$ws.Cells.Item(15, 15).Value = "Provincia del Elquí"   # O15
$ws.Cells.Item(15, 16).Value = 470   # P15

$ws.Cells.Item(16, 4).Value = 44473   # D16
$ws.Cells.Item(16, 10).Value = 500   # J16
$ws.Cells.Item(16, 11).Value = 8500   # K16
$ws.Cells.Item(16, 12).Value = 9000   # L16
$ws.Cells.Item(16, 13).Value = 8750   # M16
$ws.Cells.Item(16, 15).Value = "Provincia del Elquí"   # O16
$ws.Cells.Item(16, 16).Value = 350   # P16

$ws.Cells.Item(17, 4).Value = 44386   # D17
$ws.Cells.Item(17, 10).Value = 500   # J17
$ws.Cells.Item(17, 11).Value = 11000   # K17
$ws.Cells.Item(17, 12).Value = 12000   # L17
$ws.Cells.Item(17, 13).Value = 11500   # M17
$ws.Cells.Item(17, 15).Value = "Provincia del Elquí"   # O17
$ws.Cells.Item(17, 16).Value = 460   # P17

$ws.Cells.Item(18, 4).Value = 44690   # D18
$ws.Cells.Item(18, 10).Value = 400   # J18
$ws.Cells.Item(18, 11).Value = 17000   # K18
$ws.Cells.Item(18, 12).Value = 18000   # L18
$ws.Cells.Item(18, 13).Value = 17500   # M18
$ws.Cells.Item(18, 15).Value = "Provincia del Elquí"   # O18
$ws.Cells.Item(18, 16).Value = 700   # P18

